$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.827.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.91%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.705.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.90%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.25%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'315.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.18%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.73%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.4045"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.19%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'BinanceUSD"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'1.001"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.39%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'Polygon"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'1.472"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.73%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'OKB"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'53.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.09%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08817"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.67%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'26.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +4.66%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.517"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.47%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.985"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.32%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.00001343"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'1.716.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.65%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'95.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.95%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.75%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'20.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +5.08%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.09%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.16%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'14.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.36%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'24.821.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.90%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.366"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.46%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.886"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.28%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'23.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.18%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.206"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +18.84%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'161.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.40%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'144.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.12%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'8.200"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -7.59%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.279"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +15.46%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.905.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.82%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.08650"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.39%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.03201"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +9.42%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'7.286"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.45%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.030"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.03%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +3.72%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.8380"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +6.91%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.09484"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.79%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'10.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.74%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.00%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.479"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.25%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'17.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.56%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.715"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.59%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +3.12%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'4.219"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.55%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.379"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.04%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.25%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'140.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.65%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.08403"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +5.39%  "
$ws.Range("E51").Style = "Normal"
